# Apply updated dSF (column F) values for specific rows, per the data repull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    6  = -1
    9  = 0
    10 = -2
    12 = -1
    13 = -3
    14 = -1
    20 = 2
    21 = 1
    23 = 2
    26 = -4
    28 = -2
    29 = -7
    31 = 2
    32 = 4
    33 = -1
    40 = 2
    41 = -3
    48 = -4
    52 = -3
    53 = 0
    57 = 3
    61 = -3
    62 = 0
    63 = -2
    64 = 1
}

foreach ($row in $changes.Keys) {
    $ws.Range("F$row").Value = $changes[$row]
}
